# New updates for home page
# Update the "Actual" (column L) results on the TC_OpenBrowser sheet
# so they reflect the Action + Input1 that was actually performed.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TC_OpenBrowser")

$ws.Range("L2").Value = "OpenBrowser: Chrome"
$ws.Range("L3").Value = "GetURL: https://www.citymarketnorwalk.com/"
$ws.Range("L4").Value = "Click: City Marketplace"
